$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Vostro"/"Dell" (DKS) device-spec entries (rows 6 & 7). This
# shifts every row below up by two and causes the shared-string table to be
# recompacted when the workbook is saved.
$ws.Rows("6:7").Delete()

# Set the printable page to Letter/A4-style portrait layout (matches the
# pageSetup block written when the workbook was last saved from Excel).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the cursor/selection position recorded in the saved view state.
$ws.Range("E16").Select()
